# Add DCGAN layers and new config
# ---------------------------------------------------------------------------
# This script reproduces, via Excel COM interop, the changes described by the
# target diff:
#   * Two new shared strings ("Generator DCGAN", "100->512") used as labels
#     for a new DCGAN generator block added to the layer calculator.
#   * A new "Generator DCGAN" TransposedConv block in columns O:W, rows 63-71
#     (mirrors the existing D:L "Generator beefed-up" block just above it).
#   * A second, new CONV-based block in columns N:W, rows 75-83 (mirrors the
#     existing D:L CONV block that ends at row 75, continuing with brand new
#     rows 76-81 plus the already-existing rows 82/83).
#   * Updated sheet view (scroll position + active cell/selection).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: copy ONLY the number/cell formatting from $src to $dst without
# touching the values/formulas already present (or about to be written) in
# the destination range. Mirrors how the source workbook reuses style ids.
# NOTE: always pass plain string variables (not inline concatenation
# expressions) as arguments here.
# ---------------------------------------------------------------------------
function Copy-Format($src, $dst) {
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 1) New "Generator DCGAN" title block, row 63 (mirrors D63, style of D63).
# ---------------------------------------------------------------------------
Copy-Format "D63:L63" "O63:W63"
$ws.Range("O63").Value = "Generator DCGAN"

# ---------------------------------------------------------------------------
# 2) New header row 64 (size / kernel / s / padding / dilation), mirrors D64:L64
# ---------------------------------------------------------------------------
Copy-Format "D64:L64" "O64:W64"
$ws.Range("O64").Value = "size"
$ws.Range("P64").Value = "kernel"
$ws.Range("Q64").Value = "s"
$ws.Range("S64").Value = "padding"
$ws.Range("U64").Value = "dilation"

# ---------------------------------------------------------------------------
# 3) DCGAN TransposedConv rows 65-71 (columns N-W).
#    N = running channel count label/value, O:V = inputs, W = output formula.
# ---------------------------------------------------------------------------
Copy-Format "D65:L65" "N65:W65"
$ws.Range("N65").Value = "100→512"
$ws.Range("O65").Value = 1
$ws.Range("P65").Value = 4
$ws.Range("Q65").Value = 1
$ws.Range("S65").Value = 1
$ws.Range("U65").Value = 0
$ws.Range("W65").Formula = "=(O65-1)*Q65-2*U65+S65*(P65-1)+0+1"

$dcganRows = @(
    @{ Row = 66; N = 512; P = 4; Q = 2; S = 1; U = 1 },
    @{ Row = 67; N = 256; P = 4; Q = 2; S = 1; U = 1 },
    @{ Row = 68; N = 256; P = 4; Q = 2; S = 1; U = 1 },
    @{ Row = 69; N = 128; P = 4; Q = 2; S = 1; U = 1 },
    @{ Row = 70; N = 64;  P = 4; Q = 2; S = 1; U = 1 },
    @{ Row = 71; N = 3;   P = 4; Q = 2; S = 1; U = 1 }
)

$prev = 65
foreach ($r in $dcganRows) {
    $row = $r.Row
    $srcFmt = "D" + $row + ":L" + $row
    $dstFmt = "N" + $row + ":W" + $row
    Copy-Format $srcFmt $dstFmt
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Formula = "=W$prev"
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("S$row").Value = $r.S
    $ws.Range("U$row").Value = $r.U
    $ws.Range("W$row").Formula = "=(O$row-1)*Q$row-2*U$row+S$row*(P$row-1)+0+1"
    $prev = $row
}

# ---------------------------------------------------------------------------
# 4) New header row 75 (size / kernel / s / padding / dilation), columns O:W.
# ---------------------------------------------------------------------------
Copy-Format "D64:L64" "O75:W75"
$ws.Range("O75").Value = "size"
$ws.Range("P75").Value = "kernel"
$ws.Range("Q75").Value = "s"
$ws.Range("S75").Value = "padding"
$ws.Range("U75").Value = "dilation"

# ---------------------------------------------------------------------------
# 5) New rows 76-81 (brand new rows, columns N:W only) plus the N:W additions
#    to the already-existing rows 82 and 83.
# ---------------------------------------------------------------------------
$newConvRows = @(
    @{ Row = 76; N = 512; IsFirst = $true;  Ref = "";    P = 4; Q = 1; S = 1; U = 0 },
    @{ Row = 77; N = 512; IsFirst = $false; Ref = "W76"; P = 3; Q = 1; S = 1; U = 1 },
    @{ Row = 78; N = 512; IsFirst = $false; Ref = "W77"; P = 3; Q = 2; S = 1; U = 1 },
    @{ Row = 79; N = 256; IsFirst = $false; Ref = "W78"; P = 3; Q = 2; S = 1; U = 1 },
    @{ Row = 80; N = 256; IsFirst = $false; Ref = "W79"; P = 3; Q = 2; S = 1; U = 0 },
    @{ Row = 81; N = 128; IsFirst = $false; Ref = "W80"; P = 3; Q = 2; S = 1; U = 0 },
    @{ Row = 82; N = 64;  IsFirst = $false; Ref = "W81"; P = 3; Q = 2; S = 1; U = 0 },
    @{ Row = 83; N = 3;   IsFirst = $false; Ref = "W82"; P = 4; Q = 2; S = 1; U = 0 }
)

foreach ($r in $newConvRows) {
    $row = $r.Row
    $dstFmt = "N" + $row + ":W" + $row
    Copy-Format "N65:W65" $dstFmt
    $ws.Range("N$row").Value = $r.N
    if ($r.IsFirst) {
        $ws.Range("O$row").Value = 1
    } else {
        $ws.Range("O$row").Formula = "=" + $r.Ref
    }
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("S$row").Value = $r.S
    $ws.Range("U$row").Value = $r.U
    $ws.Range("W$row").Formula = "=(O$row-1)*Q$row-2*U$row+S$row*(P$row-1)+0+1"
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6) Sheet view: scroll position + active cell / selection.
# ---------------------------------------------------------------------------
$ws.Activate()
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("L50") } catch {}
try {
    $excel.ActiveWindow.ScrollRow = 50
    $excel.ActiveWindow.ScrollColumn = 12
} catch {}
$ws.Range("Q80").Select() | Out-Null

Write-Output "DCGAN layers + new config added"
